$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Find the next empty row below the existing log entries (row 65).
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
$row = $lastRow + 1

# Clone the formatting (style 3: centered alignment, used by every data row)
# from the previous row onto the new row before writing values into it.
$srcRange = $ws.Range($ws.Cells.Item($lastRow, 1), $ws.Cells.Item($lastRow, 8))
$dstRange = $ws.Range($ws.Cells.Item($row, 1), $ws.Cells.Item($row, 8))
$srcRange.Copy()
$dstRange.PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Append the new run-log entry.
$ws.Cells.Item($row, 1).Value = "2025-08-27 13:03:10 UTC"
$ws.Cells.Item($row, 2).Value = "2025-08-27 18:33:10 IST"
$ws.Cells.Item($row, 3).Value = "SKIPPED"
$ws.Cells.Item($row, 4).Value = "No change in PDF. Skipping download & Excel update."
$ws.Cells.Item($row, 5).Value = "https://nalcoindia.com/wp-content/uploads/2019/01/INGOT-21-08-2025.pdf"
$ws.Cells.Item($row, 6).Value = ""
$ws.Cells.Item($row, 7).Value = 0
$ws.Cells.Item($row, 8).Value = ""
